# Updates symbol list figures (Price / Volume(1h)) for Fri Jan 20 2023 run.
# Values are kept as literal text (matching the source sheet's inline-string cells),
# so each cell is first forced to Text format, written, then the format flag is cleared
# again (ClearFormats) so no residual "@"/Text number-format is left applied to the cell.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Cell = 'D2'; Value = '301.63' },
    @{ Cell = 'E2'; Value = '2.63%' },
    @{ Cell = 'D3'; Value = '32.01' },
    @{ Cell = 'E3'; Value = '2.96%' },
    @{ Cell = 'D4'; Value = '5.030' },
    @{ Cell = 'E4'; Value = '2.03%' },
    @{ Cell = 'D5'; Value = '0.07794' },
    @{ Cell = 'E5'; Value = '6.26%' },
    @{ Cell = 'D6'; Value = '2.314' },
    @{ Cell = 'E6'; Value = '1.27%' },
    @{ Cell = 'D7'; Value = '7.963' },
    @{ Cell = 'E7'; Value = '3.58%' },
    @{ Cell = 'D8'; Value = '0.9278' },
    @{ Cell = 'E8'; Value = '1.79%' },
    @{ Cell = 'D9'; Value = '0.1019' },
    @{ Cell = 'E9'; Value = '25.23%' },
    @{ Cell = 'D10'; Value = '0.1768' },
    @{ Cell = 'E10'; Value = '5.07%' },
    @{ Cell = 'D11'; Value = '0.08433' },
    @{ Cell = 'E11'; Value = '2.55%' },
    @{ Cell = 'D12'; Value = '0.03377' },
    @{ Cell = 'E12'; Value = '8.77%' },
    @{ Cell = 'D13'; Value = '0.09878' },
    @{ Cell = 'E13'; Value = '-1.84%' },
    @{ Cell = 'D14'; Value = '0.001474' },
    @{ Cell = 'E14'; Value = '-2.35%' },
    @{ Cell = 'D15'; Value = '0.005748' },
    @{ Cell = 'E15'; Value = '0.31%' },
    @{ Cell = 'D16'; Value = '3.493' },
    @{ Cell = 'E16'; Value = '0.29%' },
    @{ Cell = 'D17'; Value = '3.859' },
    @{ Cell = 'E17'; Value = '2.36%' },
    @{ Cell = 'E18'; Value = '5.35%' },
    @{ Cell = 'D19'; Value = '0.3363' },
    @{ Cell = 'E19'; Value = '1.05%' },
    @{ Cell = 'D20'; Value = '0.1344' },
    @{ Cell = 'E20'; Value = '3.05%' },
    @{ Cell = 'D21'; Value = '4.269' },
    @{ Cell = 'E21'; Value = '7.58%' },
    @{ Cell = 'E22'; Value = '-0.92%' },
    @{ Cell = 'E23'; Value = '1.42%' },
    @{ Cell = 'D24'; Value = '0.001216' },
    @{ Cell = 'E24'; Value = '0.40%' },
    @{ Cell = 'D25'; Value = '0.004382' },
    @{ Cell = 'E25'; Value = '1.00%' },
    @{ Cell = 'D26'; Value = '0.0001289' },
    @{ Cell = 'E26'; Value = '-0.86%' },
    @{ Cell = 'D27'; Value = '0.0003365' },
    @{ Cell = 'E27'; Value = '-0.87%' },
    @{ Cell = 'D39'; Value = '0.01710' },
    @{ Cell = 'E39'; Value = '6.72%' },
    @{ Cell = 'D40'; Value = '0.04739' },
    @{ Cell = 'E40'; Value = '6.76%' },
    @{ Cell = 'D41'; Value = '0.007819' },
    @{ Cell = 'E41'; Value = '7.01%' },
    @{ Cell = 'D42'; Value = '0.009749' },
    @{ Cell = 'E42'; Value = '11.57%' },
    @{ Cell = 'D43'; Value = '0.1397' },
    @{ Cell = 'E43'; Value = '5.45%' },
    @{ Cell = 'D44'; Value = '0.002061' },
    @{ Cell = 'E44'; Value = '0.09%' },
    @{ Cell = 'D45'; Value = '0.009652' },
    @{ Cell = 'E45'; Value = '4.89%' },
    @{ Cell = 'D46'; Value = '0.00006105' },
    @{ Cell = 'E46'; Value = '2.55%' },
    @{ Cell = 'D47'; Value = '0.00000000744' },
    @{ Cell = 'E47'; Value = '-0.86%' },
    @{ Cell = 'D48'; Value = '2.655' },
    @{ Cell = 'E48'; Value = '18.46%' },
    @{ Cell = 'D49'; Value = '0.001983' },
    @{ Cell = 'E49'; Value = '-31.59%' },
    @{ Cell = 'D50'; Value = '0.00002082' },
    @{ Cell = 'E50'; Value = '-0.86%' },
    @{ Cell = 'D51'; Value = '0.0001983' },
    @{ Cell = 'E51'; Value = '-0.86%' }
)

foreach ($u in $updates) {
    $r = $ws.Range($u.Cell)
    $r.NumberFormat = "@"   # force text interpretation so numeric-looking strings
                            # (e.g. "301.63") and percents (e.g. "2.63%") are not
                            # auto-converted into numbers/percentages by Excel.
    $r.Value = $u.Value
    $r.ClearFormats()       # drop the temporary Text number-format again; the source
                            # cells carry no explicit style, so we restore that state.
}
